# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '37.901.21'
$ws.Range("E2").Value = '  +1.82%  '

# Row 3
$ws.Range("D3").Value = '2.106.57'
$ws.Range("E3").Value = '  +2.36%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '''234.15'
$ws.Range("E5").Value = '  +0.57%  '

# Row 6
$ws.Range("E6").Value = '  +0.74%  '

# Row 7
$ws.Range("B7").Value = 'Solana'
$ws.Range("C7").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D7").Value = '''57.99'
$ws.Range("E7").Value = '  +1.25%  '

# Row 8
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.06%  '

# Row 9
$ws.Range("D9").Value = '''0.392'
$ws.Range("E9").Value = '  +2.18%  '

# Row 10
$ws.Range("D10").Value = '''0.0781'
$ws.Range("E10").Value = '  +2.97%  '

# Row 11
$ws.Range("E11").Value = '  +3.42%  '

# Row 12
$ws.Range("D12").Value = '2.403.01'
$ws.Range("E12").Value = '  +1.77%  '

# Row 13
$ws.Range("D13").Value = '''14.61'
$ws.Range("E13").Value = '  -1.17%  '

# Row 14
$ws.Range("D14").Value = '''21.41'
$ws.Range("E14").Value = '  +2.36%  '

# Row 15
$ws.Range("D15").Value = '''0.779'
$ws.Range("E15").Value = '  -0.51%  '

# Row 16
$ws.Range("D16").Value = '''5.28'
$ws.Range("E16").Value = '  +2.14%  '

# Row 17
$ws.Range("D17").Value = '2.100.49'
$ws.Range("E17").Value = '  +2.04%  '

# Row 18
$ws.Range("D18").Value = '37.844.42'
$ws.Range("E18").Value = '  +1.79%  '

# Row 19
$ws.Range("D19").Value = '''6.19'
$ws.Range("E19").Value = '  -2.66%  '

# Row 20
$ws.Range("D20").Value = '''71.07'
$ws.Range("E20").Value = '  +2.62%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0827'
$ws.Range("E21").Value = '  +2.04%  '

# Row 22
$ws.Range("D22").Value = '''228.21'
$ws.Range("E22").Value = '  +1.12%  '

# Row 23
$ws.Range("E23").Value = '  +0.01%  '

# Row 24
$ws.Range("E24").Value = '  -0.08%  '

# Row 25
$ws.Range("E25").Value = '  +0.05%  '

# Row 26
$ws.Range("D26").Value = '''168.19'
$ws.Range("E26").Value = '  +1.41%  '

# Row 27
$ws.Range("D27").Value = '''0.141'
$ws.Range("E27").Value = '  +10.97%  '

# Row 28
$ws.Range("D28").Value = '''9.01'
$ws.Range("E28").Value = '  +2.37%  '

# Row 29
$ws.Range("E29").Value = '  -1.75%  '

# Row 30
$ws.Range("D30").Value = '''19.56'
$ws.Range("E30").Value = '  +2.81%  '

# Row 31
$ws.Range("E31").Value = '  +0.99%  '

# Row 32
$ws.Range("E32").Value = '  +5.39%  '

# Row 33
$ws.Range("D33").Value = '''0.0632'
$ws.Range("E33").Value = '  +2.16%  '

# Row 34
$ws.Range("D34").Value = '''4.68'
$ws.Range("E34").Value = '  +1.55%  '

# Row 35
$ws.Range("E35").Value = '  +3.49%  '

# Row 36
$ws.Range("E36").Value = '  +5.80%  '

# Row 37
$ws.Range("E37").Value = '  +4.48%  '

# Row 38
$ws.Range("D38").Value = '''1.00'
$ws.Range("E38").Value = '  -0.06%  '

# Row 39
$ws.Range("E39").Value = '  -4.59%  '

# Row 40
$ws.Range("D40").Value = '''0.0993'
$ws.Range("E40").Value = '  +6.69%  '

# Row 41
$ws.Range("D41").Value = '''2.95'
$ws.Range("E41").Value = '  -0.02%  '

# Row 42
$ws.Range("D42").Value = '''97.97'
$ws.Range("E42").Value = '  +1.58%  '

# Row 43
$ws.Range("E43").Value = '  +2.24%  '

# Row 44
$ws.Range("D44").Value = '1.458.71'
$ws.Range("E44").Value = '  -1.00%  '

# Row 45
$ws.Range("D45").Value = '''1.17'
$ws.Range("E45").Value = '  -0.53%  '

# Row 46
$ws.Range("E46").Value = '  +4.58%  '

# Row 47
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '''15.77'
$ws.Range("E47").Value = '  +4.44%  '

# Row 48
$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D48").Value = '''4.08'
$ws.Range("E48").Value = '  -7.49%  '

# Row 49
$ws.Range("D49").Value = '''7.34'
$ws.Range("E49").Value = '  +2.49%  '

# Row 50
$ws.Range("E50").Value = '  +2.60%  '

# Row 51
$ws.Range("D51").Value = '2.300.56'
$ws.Range("E51").Value = '  +2.30%  '
